# Apply the "Sub-Count" -> "Sub-County" header rename and move the
# worksheet's active selection/scroll position over to the renamed
# column (M1), matching the author's "enhancement of public dashboard
# and reporting" edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the mis-spelled column header (shared string used by M1) ---
$ws.Range("M1").Value = "Sub-County"

# --- 2. Reflect the scrolled/selected view from the diff: the user
#        scrolled right so column F is the left-most visible column and
#        selected M1 (was topLeftCell A1 / selection H7). Attempt the
#        scroll-position update via the Window object, then select M1.
$win = $excel.ActiveWindow
try {
    $win.ScrollColumn = 6
    $win.ScrollRow = 1
} catch {
    # Scroll-position isn't always settable in every host; ignore.
}

$ws.Range("M1").Select()
